# Generate Report for Handoff
# Update status text and timestamps across the three worksheets of the
# localization-status workbook.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) go from
# "Handed back: in sync with en-US" to "Ready for handoff", and the
# "Latest HO Xliff Generate Date" column (G2) gets a refreshed timestamp.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-24 15:09:10"

# zh-cn sheet: Status (C2) matches the same "Ready for handoff" text, and
# Latest Handoff Datetime (H2) timestamp refreshed.
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-24 15:08:57"

# de-de sheet: Status (C2) matches the same "Ready for handoff" text, and
# Latest Handoff Datetime (H2) matches the refreshed Overview timestamp.
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-24 15:09:10"

# Column widths shrink to fit the new (shorter) "Ready for handoff" text
# (was sized for "Handed back: in sync with en-US").
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333336
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333336

$wb.Save()
